$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.973.67"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.692.76"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "650.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "4.319.55"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "3.679.70"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "69.966.94"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "3.839.97"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000127"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "3.689.92"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "180.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.933"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -4.43%  "
